# edit.ps1 - applies the diff:
#  1) Slide 1 (title): "обработки сигнала" -> "обработки сигналов"
#  2) Slide 2, shape 22: reposition/resize textbox + same word fix in its text
#  3) Slide 7, shape 17: width tweak only (1 EMU)
#  4) Slide 9, shape 6: left offset tweak only (1 EMU)

$p = $ppt.ActivePresentation

# --- 1) Slide 1 title: "обработки сигнала" -> "обработки сигналов" ---
# The title text is split across three runs joined by <a:br/> soft line
# breaks ("Программный модуль" / "обработки сигнала" / "стандарта DMR").
# Replacing the whole TextRange.Text would collapse the <a:br/> breaks into
# separate paragraphs, so instead only the affected run's characters are
# replaced via the Characters() sub-range.
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Characters(20, 17).Text = "обработки сигналов"

# --- 2) Slide 2, shape 22: move/resize + text fix ---
# Unlike the title on slide 1, this shape's text is one single run with no
# <a:br/> breaks, so replacing the whole TextRange.Text keeps it as a single
# run (matching the diff) instead of fragmenting it into multiple runs.
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(22)
$sh2.TextFrame.TextRange.Text = "Реализация программного модуля обработки сигналов стандарта DMR"
$sh2.Left = 41.04191207885742
$sh2.Top = 98.40083312988281
$sh2.Width = 561.4547119140625
$sh2.Height = 28.828267716535432

# --- 3) Slide 7, shape 17: width changes by 1 EMU ---
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(17)
$sh7.Width = 265.2879333496094

# --- 4) Slide 9, shape 6: left offset changes by 1 EMU ---
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(6)
$sh9.Left = 306.5577697753906
